$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = -11.854
$ws.Range("C13").Value = -12.201
$ws.Range("C16").Value = -11.983
$ws.Range("C18").Value = -12.01
$ws.Range("C20").Value = -12.32
